# Update crypto price ("Price") and volume-change ("Volume(1h)") columns
# with the latest scraped values (GitHub Actions refresh).
# Note: Price values are prefixed with a literal leading apostrophe so
# Excel stores them as text (preserving formats like "0.600", "41.60",
# and multi-dot big numbers like "67.283.85") instead of coercing them
# into numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.283.85"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "'3.483.32"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'593.76"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'178.23"
$ws.Range("E6").Value = "  +3.47%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").Value = "'3.486.49"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("E10").Value = "  +4.57%  "
$ws.Range("E11").Value = "  -2.30%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "'4.085.65"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "'31.88"
$ws.Range("E14").Value = "  +9.54%  "
$ws.Range("D15").Value = "'0.135"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "'67.303.08"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "'3.481.98"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "'388.40"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  -3.87%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("D33").Value = "'2.05"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "'23.51"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "'1.60"
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("D38").Value = "'164.62"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("E41").Value = "  +7.96%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "'2.833.22"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'27.02"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").Value = "'26.12"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D48").Value = "'41.60"
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("E51").Value = "  -2.40%  "
